$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.666.74"
$ws.Range("E2").Value = "  +5.89%  "
$ws.Range("D3").Value = "2.045.80"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.39"
$ws.Range("E5").Value = "  +4.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.649"
$ws.Range("E6").Value = "  +3.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "65.44"
$ws.Range("E7").Value = "  +15.76%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "59.47"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.376"
$ws.Range("E10").Value = "  +4.65%  "
$ws.Range("E11").Value = "  +3.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.902"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.06"
$ws.Range("E14").Value = "  +6.03%  "
$ws.Range("D15").Value = "2.342.47"
$ws.Range("E15").Value = "  +3.26%  "
$ws.Range("E16").Value = "  +6.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.36"
$ws.Range("E17").Value = "  +18.48%  "
$ws.Range("D18").Value = "2.034.55"
$ws.Range("E18").Value = "  +2.51%  "
$ws.Range("D19").Value = "37.602.59"
$ws.Range("E19").Value = "  +6.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.49"
$ws.Range("E20").Value = "  +4.72%  "
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.32"
$ws.Range("E22").Value = "  +5.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.26"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("E24").Value = "  +17.36%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("E26").Value = "  +4.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  +5.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.79"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.88"
$ws.Range("E29").Value = "  +2.18%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.114"
$ws.Range("E30").Value = "  +27.38%  "
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.20"
$ws.Range("E32").Value = "  +8.65%  "
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.70"
$ws.Range("E34").Value = "  +10.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0612"
$ws.Range("E35").Value = "  +4.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  +6.60%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +25.16%  "
$ws.Range("E39").Value = "  +2.20%  "
$ws.Range("E40").Value = "  +15.88%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.80"
$ws.Range("E41").Value = "  +24.27%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.23"
$ws.Range("E42").Value = "  +3.83%  "
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0218"
$ws.Range("E44").Value = "  +4.58%  "
$ws.Range("E45").Value = "  +5.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.06"
$ws.Range("E46").Value = "  +8.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.81"
$ws.Range("E47").Value = "  +9.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "94.84"
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("D49").Value = "1.421.14"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("E50").Value = "  +2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.38"
$ws.Range("E51").Value = "  +3.85%  "
